$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks like a plain number (e.g. "244.57") must be
# forced to Text format first, otherwise Excel auto-converts the typed
# string into a numeric value and mangles things like leading/trailing zeros
# (e.g. "21.30" -> 21.3, "0.0211" -> 2.11E-02). The NumberFormat is reset back
# to the sheet default afterwards via Style="Normal" so no stray formatting
# is left behind on the cell.
$textCells = 'D5','D6','D7','D10','D12','D13','D15','D19','D21','D23','D24','D25','D27','D28','D29','D30','D35','D38','D39','D43','D45','D46','D48','D51'
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '36.535.55'
$ws.Range('E2').Value = '  +0.31%  '
$ws.Range('D3').Value = '1.957.24'
$ws.Range('E3').Value = '  +1.08%  '
$ws.Range('E4').Value = '  -0.09%  '
$ws.Range('D5').Value = '244.57'
$ws.Range('E5').Value = '  +1.15%  '
$ws.Range('D6').Value = '0.615'
$ws.Range('E6').Value = '  +0.84%  '
$ws.Range('D7').Value = '58.48'
$ws.Range('E7').Value = '  +2.85%  '
$ws.Range('E8').Value = '  -0.09%  '
$ws.Range('E9').Value = '  +4.17%  '
$ws.Range('D10').Value = '0.0793'
$ws.Range('E10').Value = '  -5.96%  '
$ws.Range('E11').Value = '  -0.85%  '
$ws.Range('D12').Value = '14.15'
$ws.Range('E12').Value = '  +5.55%  '
$ws.Range('D13').Value = '0.837'
$ws.Range('E13').Value = '  +3.57%  '
$ws.Range('D14').Value = '2.244.46'
$ws.Range('E14').Value = '  +0.99%  '
$ws.Range('D15').Value = '21.30'
$ws.Range('E15').Value = '  +0.79%  '
$ws.Range('E16').Value = '  +2.47%  '
$ws.Range('D17').Value = '1.956.53'
$ws.Range('E17').Value = '  +1.00%  '
$ws.Range('D18').Value = '36.512.44'
$ws.Range('E18').Value = '  +0.48%  '
$ws.Range('D19').Value = '69.77'
$ws.Range('E19').Value = '  +0.79%  '
$ws.Range('D20').Value = '0.0₃0848'
$ws.Range('E20').Value = '  -1.95%  '
$ws.Range('D21').Value = '229.34'
$ws.Range('E21').Value = '  +0.55%  '
$ws.Range('E22').Value = '  +1.42%  '
$ws.Range('D23').Value = '0.999'
$ws.Range('E23').Value = '  -0.08%  '
$ws.Range('D24').Value = '2.47'
$ws.Range('E24').Value = '  +5.07%  '
$ws.Range('D25').Value = '2.36'
$ws.Range('E25').Value = '  +3.78%  '
$ws.Range('E26').Value = '  +8.04%  '
$ws.Range('D27').Value = '9.14'
$ws.Range('E27').Value = '  -1.13%  '
$ws.Range('D28').Value = '160.45'
$ws.Range('E28').Value = '  -0.21%  '
$ws.Range('D29').Value = '19.40'
$ws.Range('E29').Value = '  +1.07%  '
$ws.Range('D30').Value = '0.119'
$ws.Range('E30').Value = '  +1.65%  '
$ws.Range('E31').Value = '  +6.26%  '
$ws.Range('E32').Value = '  +3.21%  '
$ws.Range('E33').Value = '  -2.79%  '
$ws.Range('E34').Value = '  +5.63%  '
$ws.Range('D35').Value = '3.46'
$ws.Range('E35').Value = '  +15.73%  '
$ws.Range('E36').Value = '  +7.30%  '
$ws.Range('E37').Value = '  -0.02%  '
$ws.Range('D38').Value = '1.76'
$ws.Range('E38').Value = '  -1.38%  '
$ws.Range('D39').Value = '5.40'
$ws.Range('E39').Value = '  -12.07%  '
$ws.Range('E40').Value = '  +0.40%  '
$ws.Range('E41').Value = '  +1.52%  '
$ws.Range('E42').Value = '  +0.75%  '
$ws.Range('D43').Value = '0.0211'
$ws.Range('E43').Value = '  +0.93%  '
$ws.Range('D44').Value = '1.375.48'
$ws.Range('E44').Value = '  +2.86%  '
$ws.Range('D45').Value = '15.71'
$ws.Range('E45').Value = '  +1.22%  '
$ws.Range('D46').Value = '88.07'
$ws.Range('E46').Value = '  +1.35%  '
$ws.Range('E47').Value = '  +0.36%  '
$ws.Range('D48').Value = '7.13'
$ws.Range('E48').Value = '  +0.48%  '
$ws.Range('E49').Value = '  +0.57%  '
$ws.Range('D50').Value = '2.135.06'
$ws.Range('E50').Value = '  +0.93%  '
$ws.Range('B51').Value = 'FTXToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D51').Value = '3.50'
$ws.Range('E51').Value = '  +18.93%  '

foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
